$wb = $excel.ActiveWorkbook

# Add a new worksheet placed right after the last sheet (SimpleDemNest),
# which makes it the active/selected tab, matching the target workbook state.
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "TwoxTwowAuxDem"

# Populate cells in the same order the source values were entered so that
# newly-introduced shared strings land in the expected sequence.
$newSheet.Range("A11").Value = "'SXX.L"
$newSheet.Range("A12").Value = "'SYY.L"
$newSheet.Range("A13").Value = "'SWW.L"
$newSheet.Range("A14").Value = "'DKX.L"
$newSheet.Range("A15").Value = "'DLX.L"
$newSheet.Range("A16").Value = "'DLY.L"
$newSheet.Range("A17").Value = "'DKY.L"
$newSheet.Range("A18").Value = "'DXW.L"
$newSheet.Range("A19").Value = "'DYW.L"
$newSheet.Range("A21").Value = "'CWCONS.L"
$newSheet.Range("C1").Value = "'UnEmp=.1"
$newSheet.Range("D1").Value = "'TKX=0.25"
$newSheet.Range("F1").Value = "'TY=0.5"
$newSheet.Range("E1").Value = "'&TY=.5"
$newSheet.Range("B1").Value = "'benchmark"
$newSheet.Range("A2").Value = "'X.L"
$newSheet.Range("B2").Value = 1
$newSheet.Range("C2").Value = 1.0000000000000002
$newSheet.Range("D2").Value = 1.0371252977019936
$newSheet.Range("E2").Value = 0.98447088030257324
$newSheet.Range("F2").Value = 1.0069367675479344
$newSheet.Range("A3").Value = "'Y.L"
$newSheet.Range("B3").Value = 1
$newSheet.Range("C3").Value = 1.0000000000000004
$newSheet.Range("D3").Value = 1.053409075564238
$newSheet.Range("E3").Value = 0.44661553192056858
$newSheet.Range("F3").Value = 0.41177381270347702
$newSheet.Range("A4").Value = "'W.L"
$newSheet.Range("B4").Value = 1
$newSheet.Range("C4").Value = 1.0000000000000002
$newSheet.Range("D4").Value = 1.0452354763862961
$newSheet.Range("E4").Value = 0.66308369432539305
$newSheet.Range("F4").Value = 0.64391784561311827
$newSheet.Range("A5").Value = "'PX.L"
$newSheet.Range("B5").Value = 1
$newSheet.Range("C5").Value = 1.0499999999722391
$newSheet.Range("D5").Value = 1
$newSheet.Range("E5").Value = 1
$newSheet.Range("F5").Value = 1
$newSheet.Range("A6").Value = "'PY.L"
$newSheet.Range("B6").Value = 1
$newSheet.Range("C6").Value = 1.0499999999722391
$newSheet.Range("D6").Value = 0.9845418287728368
$newSheet.Range("E6").Value = 2.2042916262845966
$newSheet.Range("F6").Value = 2.4453637808572513
$newSheet.Range("A7").Value = "'PW.L"
$newSheet.Range("B7").Value = 1
$newSheet.Range("C7").Value = 1.0499999999722391
$newSheet.Range("D7").Value = 0.99224081188632707
$newSheet.Range("E7").Value = 1.4846856995121664
$newSheet.Range("F7").Value = 1.5637658970763419
$newSheet.Range("A8").Value = "'PL.L"
$newSheet.Range("B8").Value = 1
$newSheet.Range("C8").Value = 1.0499999999722391
$newSheet.Range("D8").Value = 0.99224081188632707
$newSheet.Range("E8").Value = 1.4846856995121667
$newSheet.Range("F8").Value = 1.5637658970763419
$newSheet.Range("A9").Value = "'PK.L"
$newSheet.Range("B9").Value = 1
$newSheet.Range("C9").Value = 1.0499999999722394
$newSheet.Range("D9").Value = 0.97310521759683544
$newSheet.Range("E9").Value = 0.704929765951823
$newSheet.Range("F9").Value = 0.84532963183954124
$newSheet.Range("A10").Value = "'U.L"
$newSheet.Range("B10").Value = 0.2
$newSheet.Range("C10").Value = 0.19999999999999957
$newSheet.Range("D10").Value = 0.10830961355206511
$newSheet.Range("E10").Value = 0.70637102840442134
$newSheet.Range("F10").Value = 0.72063893707626925
$newSheet.Range("B11").Value = 100
$newSheet.Range("C11").Value = 100
$newSheet.Range("D11").Value = 100
$newSheet.Range("E11").Value = 100
$newSheet.Range("F11").Value = 100
$newSheet.Range("B12").Value = 100
$newSheet.Range("C12").Value = 100
$newSheet.Range("D12").Value = 100
$newSheet.Range("E12").Value = 100
$newSheet.Range("F12").Value = 100
$newSheet.Range("B13").Value = 200
$newSheet.Range("C13").Value = 200
$newSheet.Range("D13").Value = 200
$newSheet.Range("E13").Value = 200
$newSheet.Range("F13").Value = 200
$newSheet.Range("B14").Value = 50
$newSheet.Range("C14").Value = 50
$newSheet.Range("D14").Value = 45.67280458551334
$newSheet.Range("E14").Value = 63.048046154598161
$newSheet.Range("F14").Value = 65.720582199051307
$newSheet.Range("B15").Value = 40
$newSheet.Range("C15").Value = 40.000000000000007
$newSheet.Range("D15").Value = 44.79199395149989
$newSheet.Range("E15").Value = 29.935254602427996
$newSheet.Range("F15").Value = 28.421418146269836
$newSheet.Range("B16").Value = 60
$newSheet.Range("C16").Value = 60.000000000000014
$newSheet.Range("D16").Value = 59.53444871317862
$newSheet.Range("E16").Value = 44.540570960010214
$newSheet.Range("F16").Value = 46.91297691147642
$newSheet.Range("B17").Value = 40
$newSheet.Range("C17").Value = 40
$newSheet.Range("D17").Value = 40.470107896625159
$newSheet.Range("E17").Value = 62.539326087466797
$newSheet.Range("F17").Value = 57.855863291935997
$newSheet.Range("B18").Value = 100
$newSheet.Range("C18").Value = 100
$newSheet.Range("D18").Value = 99.224081188632667
$newSheet.Range("E18").Value = 148.46856994948786
$newSheet.Range("F18").Value = 156.37658970757903
$newSheet.Range("B19").Value = 100
$newSheet.Range("C19").Value = 100
$newSheet.Range("D19").Value = 100.78198639087648
$newSheet.Range("E19").Value = 67.354322894079274
$newSheet.Range("F19").Value = 63.948190830224604
$newSheet.Range("A20").Value = "'CONS.L"
$newSheet.Range("B20").Value = 200
$newSheet.Range("C20").Value = 210
$newSheet.Range("D20").Value = 207.42505954038575
$newSheet.Range("E20").Value = 196.89417562455813
$newSheet.Range("F20").Value = 201.38735349872735
$newSheet.Range("B21").Value = 200
$newSheet.Range("C21").Value = 200.00000000000006
$newSheet.Range("D21").Value = 209.0470952772589
$newSheet.Range("E21").Value = 132.61673880825623
$newSheet.Range("F21").Value = 128.78356912325972
$newSheet.Range("A22").Value = "'CWI.L"
$newSheet.Range("B22").Value = 1
$newSheet.Range("C22").Value = 1.0000000000000002
$newSheet.Range("D22").Value = 1.0452354763862945
$newSheet.Range("E22").Value = 0.66308369404128109
$newSheet.Range("F22").Value = 0.64391784561629861
